$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# --- Swap/update existing rows (13 pairs) ---
$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 6230048
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 44963.5625
$arr[0,4] = 'Al Talaba'
$arr[0,5] = 'Al Karkh'
$arr[0,6] = 2
$arr[0,7] = 0
$arr[0,8] = 'H'
$arr[0,9] = 2.2
$arr[0,10] = 3
$arr[0,11] = 3.1
$arr[0,12] = 2.1
$arr[0,13] = 3
$arr[0,14] = 3.25
$arr[0,15] = -0.25
$arr[0,16] = 1.875
$arr[0,17] = 1.925
$arr[0,18] = 2
$arr[0,19] = 2.025
$arr[0,20] = 1.775
$arr[0,21] = 1.1
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 0.875
$arr[0,25] = -1
$arr[0,26] = 0
$arr[0,27] = 0
$ws.Range($ws.Cells.Item(14,2), $ws.Cells.Item(14,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 6230049
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 44963.5625
$arr[0,4] = 'Al Quwa Al Jawiya'
$arr[0,5] = 'Naft AlBasra'
$arr[0,6] = 0
$arr[0,7] = 1
$arr[0,8] = 'A'
$arr[0,9] = 1.727
$arr[0,10] = 3.4
$arr[0,11] = 4.2
$arr[0,12] = 1.7
$arr[0,13] = 3.4
$arr[0,14] = 4.333
$arr[0,15] = -0.75
$arr[0,16] = 1.975
$arr[0,17] = 1.825
$arr[0,18] = 2
$arr[0,19] = 2.025
$arr[0,20] = 1.775
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 3.333
$arr[0,24] = -1
$arr[0,25] = 0.825
$arr[0,26] = -1
$arr[0,27] = 0.7749999999999999
$ws.Range($ws.Cells.Item(15,2), $ws.Cells.Item(15,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 6315121
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 44981.33333333334
$arr[0,4] = 'Al Sinaah'
$arr[0,5] = 'Naft AlBasra'
$arr[0,6] = 1
$arr[0,7] = 0
$arr[0,8] = 'H'
$arr[0,9] = 3.4
$arr[0,10] = 2.5
$arr[0,11] = 2.4
$arr[0,12] = 3.25
$arr[0,13] = 2.55
$arr[0,14] = 2.45
$arr[0,15] = 0.25
$arr[0,16] = 1.725
$arr[0,17] = 2.075
$arr[0,18] = 1.75
$arr[0,19] = 1.8
$arr[0,20] = 2
$arr[0,21] = 2.25
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 0.7250000000000001
$arr[0,25] = -1
$arr[0,26] = -1
$arr[0,27] = 1
$ws.Range($ws.Cells.Item(25,2), $ws.Cells.Item(25,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 6315120
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 44981.33333333334
$arr[0,4] = 'Al Karkh'
$arr[0,5] = 'Al Zawraa'
$arr[0,6] = 2
$arr[0,7] = 2
$arr[0,8] = 'D'
$arr[0,9] = 2.6
$arr[0,10] = 2.6
$arr[0,11] = 2.9
$arr[0,12] = 2.9
$arr[0,13] = 2.5
$arr[0,14] = 2.7
$arr[0,15] = 0
$arr[0,16] = 2
$arr[0,17] = 1.8
$arr[0,18] = 2
$arr[0,19] = 2
$arr[0,20] = 1.8
$arr[0,21] = -1
$arr[0,22] = 1.5
$arr[0,23] = -1
$arr[0,24] = 0
$arr[0,25] = 0
$arr[0,26] = 1
$arr[0,27] = -1
$ws.Range($ws.Cells.Item(26,2), $ws.Cells.Item(26,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 6410621
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45003.375
$arr[0,4] = 'Al Sinaah'
$arr[0,5] = 'Zakho'
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 'D'
$arr[0,9] = 2.8
$arr[0,10] = 2.8
$arr[0,11] = 2.5
$arr[0,12] = 2.7
$arr[0,13] = 2.9
$arr[0,14] = 2.5
$arr[0,15] = 0
$arr[0,16] = 2
$arr[0,17] = 1.8
$arr[0,18] = 1.75
$arr[0,19] = 1.725
$arr[0,20] = 1.975
$arr[0,21] = -1
$arr[0,22] = 1.9
$arr[0,23] = -1
$arr[0,24] = 0
$arr[0,25] = 0
$arr[0,26] = -1
$arr[0,27] = 0.9750000000000001
$ws.Range($ws.Cells.Item(39,2), $ws.Cells.Item(39,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 6410619
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45003.375
$arr[0,4] = 'Al Hudod'
$arr[0,5] = 'Al Naft SC'
$arr[0,6] = 4
$arr[0,7] = 1
$arr[0,8] = 'H'
$arr[0,9] = 2.5
$arr[0,10] = 2.6
$arr[0,11] = 3
$arr[0,12] = 2.5
$arr[0,13] = 2.6
$arr[0,14] = 3
$arr[0,15] = -0.25
$arr[0,16] = 2
$arr[0,17] = 1.8
$arr[0,18] = 1.75
$arr[0,19] = 1.9
$arr[0,20] = 1.9
$arr[0,21] = 1.5
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 1
$arr[0,25] = -1
$arr[0,26] = 0.8999999999999999
$arr[0,27] = -1
$ws.Range($ws.Cells.Item(40,2), $ws.Cells.Item(40,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 6497757
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45024.39583333334
$arr[0,4] = 'Al Sinaah'
$arr[0,5] = 'Al Najaf'
$arr[0,6] = 0
$arr[0,7] = 1
$arr[0,8] = 'A'
$arr[0,9] = 3.5
$arr[0,10] = 2.8
$arr[0,11] = 2.1
$arr[0,12] = 3.5
$arr[0,13] = 2.8
$arr[0,14] = 2.1
$arr[0,15] = 0.25
$arr[0,16] = 1.975
$arr[0,17] = 1.825
$arr[0,18] = 1.75
$arr[0,19] = 1.9
$arr[0,20] = 1.9
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 1.1
$arr[0,24] = -1
$arr[0,25] = 0.825
$arr[0,26] = -1
$arr[0,27] = 0.8999999999999999
$ws.Range($ws.Cells.Item(45,2), $ws.Cells.Item(45,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 6497756
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45024.39583333334
$arr[0,4] = 'Al Karkh'
$arr[0,5] = 'Al Qasim SC'
$arr[0,6] = 2
$arr[0,7] = 1
$arr[0,8] = 'H'
$arr[0,9] = 2
$arr[0,10] = 2.875
$arr[0,11] = 3.75
$arr[0,12] = 2.15
$arr[0,13] = 2.75
$arr[0,14] = 3.4
$arr[0,15] = -0.25
$arr[0,16] = 1.9
$arr[0,17] = 1.9
$arr[0,18] = 1.75
$arr[0,19] = 1.875
$arr[0,20] = 1.925
$arr[0,21] = 1.15
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 0.8999999999999999
$arr[0,25] = -1
$arr[0,26] = 0.875
$arr[0,27] = -1
$ws.Range($ws.Cells.Item(46,2), $ws.Cells.Item(46,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 6707962
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45075.44791666666
$arr[0,4] = 'Newroz SC'
$arr[0,5] = 'Duhok'
$arr[0,6] = 2
$arr[0,7] = 1
$arr[0,8] = 'H'
$arr[0,9] = 2.2
$arr[0,10] = 2.9
$arr[0,11] = 3.2
$arr[0,12] = 2.25
$arr[0,13] = 2.75
$arr[0,14] = 3.25
$arr[0,15] = -0.25
$arr[0,16] = 1.975
$arr[0,17] = 1.825
$arr[0,18] = 1.75
$arr[0,19] = 1.825
$arr[0,20] = 1.975
$arr[0,21] = 1.25
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 0.9750000000000001
$arr[0,25] = -1
$arr[0,26] = 0.825
$arr[0,27] = -1
$ws.Range($ws.Cells.Item(80,2), $ws.Cells.Item(80,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 6704935
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45075.44791666666
$arr[0,4] = 'Al Kahrabaa'
$arr[0,5] = 'Karbalaa FC'
$arr[0,6] = 0
$arr[0,7] = 1
$arr[0,8] = 'A'
$arr[0,9] = 1.8
$arr[0,10] = 3
$arr[0,11] = 4.5
$arr[0,12] = 1.666
$arr[0,13] = 3.1
$arr[0,14] = 5
$arr[0,15] = -0.75
$arr[0,16] = 1.95
$arr[0,17] = 1.85
$arr[0,18] = 2
$arr[0,19] = 2
$arr[0,20] = 1.8
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 4
$arr[0,24] = -1
$arr[0,25] = 0.8500000000000001
$arr[0,26] = -1
$arr[0,27] = 0.8
$ws.Range($ws.Cells.Item(81,2), $ws.Cells.Item(81,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 6862618
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45112.45833333334
$arr[0,4] = 'Naft AlWasat'
$arr[0,5] = 'Newroz SC'
$arr[0,6] = 0
$arr[0,7] = 2
$arr[0,8] = 'A'
$arr[0,9] = 3
$arr[0,10] = 2.4
$arr[0,11] = 2.75
$arr[0,12] = 2.05
$arr[0,13] = 2.625
$arr[0,14] = 4
$arr[0,15] = -0.25
$arr[0,16] = 1.775
$arr[0,17] = 2.025
$arr[0,18] = 2
$arr[0,19] = 1.95
$arr[0,20] = 1.85
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 3
$arr[0,24] = -1
$arr[0,25] = 1.025
$arr[0,26] = 0
$arr[0,27] = 0
$ws.Range($ws.Cells.Item(105,2), $ws.Cells.Item(105,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 6862617
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45112.45833333334
$arr[0,4] = 'Al Naft SC'
$arr[0,5] = 'Al Najaf'
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 'D'
$arr[0,9] = 3.6
$arr[0,10] = 2.8
$arr[0,11] = 2.1
$arr[0,12] = 1.95
$arr[0,13] = 2.8
$arr[0,14] = 4.1
$arr[0,15] = -0.5
$arr[0,16] = 2
$arr[0,17] = 1.8
$arr[0,18] = 2.25
$arr[0,19] = 2.025
$arr[0,20] = 1.775
$arr[0,21] = -1
$arr[0,22] = 1.8
$arr[0,23] = -1
$arr[0,24] = -1
$arr[0,25] = 0.8
$arr[0,26] = -1
$arr[0,27] = 0.7749999999999999
$ws.Range($ws.Cells.Item(106,2), $ws.Cells.Item(106,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 7407173
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45229.64583333334
$arr[0,4] = 'Al Shorta SC'
$arr[0,5] = 'Karbalaa FC'
$arr[0,6] = 2
$arr[0,7] = 1
$arr[0,8] = 'H'
$arr[0,9] = 1.333
$arr[0,10] = 4
$arr[0,11] = 8.5
$arr[0,12] = 1.333
$arr[0,13] = 4
$arr[0,14] = 8.5
$arr[0,15] = -1.25
$arr[0,16] = 1.825
$arr[0,17] = 1.975
$arr[0,18] = 2.25
$arr[0,19] = 1.875
$arr[0,20] = 1.925
$arr[0,21] = 0.333
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = -0.5
$arr[0,25] = 0.4875
$arr[0,26] = 0.875
$arr[0,27] = -1
$ws.Range($ws.Cells.Item(124,2), $ws.Cells.Item(124,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 7406927
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45229.64583333334
$arr[0,4] = 'Naft Maysan'
$arr[0,5] = 'Al Zawraa'
$arr[0,6] = 1
$arr[0,7] = 1
$arr[0,8] = 'D'
$arr[0,9] = 3.4
$arr[0,10] = 2.8
$arr[0,11] = 2.15
$arr[0,12] = 3.4
$arr[0,13] = 2.8
$arr[0,14] = 2.15
$arr[0,15] = 0.25
$arr[0,16] = 1.9
$arr[0,17] = 1.9
$arr[0,18] = 2
$arr[0,19] = 2
$arr[0,20] = 1.8
$arr[0,21] = -1
$arr[0,22] = 1.8
$arr[0,23] = -1
$arr[0,24] = 0.45
$arr[0,25] = -0.5
$arr[0,26] = 0
$arr[0,27] = 0
$ws.Range($ws.Cells.Item(125,2), $ws.Cells.Item(125,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 7555939
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45265.35416666666
$arr[0,4] = 'Amanat Baghdad'
$arr[0,5] = 'Al Shorta SC'
$arr[0,6] = 1
$arr[0,7] = 2
$arr[0,8] = 'A'
$arr[0,9] = 6.5
$arr[0,10] = 3.75
$arr[0,11] = 1.444
$arr[0,12] = 7.5
$arr[0,13] = 3.8
$arr[0,14] = 1.4
$arr[0,15] = 1.25
$arr[0,16] = 1.825
$arr[0,17] = 1.975
$arr[0,18] = 2.25
$arr[0,19] = 1.875
$arr[0,20] = 1.925
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.3999999999999999
$arr[0,24] = 0.4125
$arr[0,25] = -0.5
$arr[0,26] = 0.875
$arr[0,27] = -1
$ws.Range($ws.Cells.Item(152,2), $ws.Cells.Item(152,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 7555937
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45265.35416666666
$arr[0,4] = 'Karbalaa FC'
$arr[0,5] = 'Al Naft SC'
$arr[0,6] = 0
$arr[0,7] = 2
$arr[0,8] = 'A'
$arr[0,9] = 2.1
$arr[0,10] = 2.875
$arr[0,11] = 3.4
$arr[0,12] = 2.1
$arr[0,13] = 2.875
$arr[0,14] = 3.4
$arr[0,15] = -0.25
$arr[0,16] = 1.85
$arr[0,17] = 1.95
$arr[0,18] = 2
$arr[0,19] = 1.8
$arr[0,20] = 2
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 2.4
$arr[0,24] = -1
$arr[0,25] = 0.95
$arr[0,26] = 0
$arr[0,27] = 0
$ws.Range($ws.Cells.Item(153,2), $ws.Cells.Item(153,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 7565124
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45268.35416666666
$arr[0,4] = 'Al Naft SC'
$arr[0,5] = 'Amanat Baghdad'
$arr[0,6] = 0
$arr[0,7] = 1
$arr[0,8] = 'A'
$arr[0,9] = 1.5
$arr[0,10] = 3.6
$arr[0,11] = 6
$arr[0,12] = 1.727
$arr[0,13] = 3.5
$arr[0,14] = 4.2
$arr[0,15] = -0.75
$arr[0,16] = 2
$arr[0,17] = 1.8
$arr[0,18] = 2
$arr[0,19] = 1.875
$arr[0,20] = 1.925
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 3.2
$arr[0,24] = -1
$arr[0,25] = 0.8
$arr[0,26] = -1
$arr[0,27] = 0.925
$ws.Range($ws.Cells.Item(156,2), $ws.Cells.Item(156,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 7565123
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45268.35416666666
$arr[0,4] = 'Al Karkh'
$arr[0,5] = 'Karbalaa FC'
$arr[0,6] = 1
$arr[0,7] = 2
$arr[0,8] = 'A'
$arr[0,9] = 1.727
$arr[0,10] = 3.4
$arr[0,11] = 4.2
$arr[0,12] = 1.833
$arr[0,13] = 3.3
$arr[0,14] = 3.75
$arr[0,15] = -0.5
$arr[0,16] = 1.875
$arr[0,17] = 1.925
$arr[0,18] = 2
$arr[0,19] = 1.925
$arr[0,20] = 1.875
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 2.75
$arr[0,24] = -1
$arr[0,25] = 0.925
$arr[0,26] = 0.925
$arr[0,27] = -1
$ws.Range($ws.Cells.Item(157,2), $ws.Cells.Item(157,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 7592803
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45275.35416666666
$arr[0,4] = 'Naft AlBasra'
$arr[0,5] = 'Newroz SC'
$arr[0,6] = 0
$arr[0,7] = 1
$arr[0,8] = 'A'
$arr[0,9] = 2.3
$arr[0,10] = 2.875
$arr[0,11] = 3
$arr[0,12] = 2.3
$arr[0,13] = 2.875
$arr[0,14] = 3
$arr[0,15] = -0.25
$arr[0,16] = 2.025
$arr[0,17] = 1.775
$arr[0,18] = 2
$arr[0,19] = 2
$arr[0,20] = 1.8
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 2
$arr[0,24] = -1
$arr[0,25] = 0.7749999999999999
$arr[0,26] = -1
$arr[0,27] = 0.8
$ws.Range($ws.Cells.Item(161,2), $ws.Cells.Item(161,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 7590595
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45275.35416666666
$arr[0,4] = 'Al Hudod'
$arr[0,5] = 'Al Shorta SC'
$arr[0,6] = 0
$arr[0,7] = 4
$arr[0,8] = 'A'
$arr[0,9] = 3.75
$arr[0,10] = 2.9
$arr[0,11] = 2
$arr[0,12] = 4.5
$arr[0,13] = 3
$arr[0,14] = 1.75
$arr[0,15] = 0.75
$arr[0,16] = 1.75
$arr[0,17] = 1.95
$arr[0,18] = 2.5
$arr[0,19] = 2
$arr[0,20] = 1.8
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.75
$arr[0,24] = -1
$arr[0,25] = 0.95
$arr[0,26] = 1
$arr[0,27] = -1
$ws.Range($ws.Cells.Item(162,2), $ws.Cells.Item(162,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 7618726
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45283.35416666666
$arr[0,4] = 'Al Hudod'
$arr[0,5] = 'Duhok'
$arr[0,6] = 0
$arr[0,7] = 1
$arr[0,8] = 'A'
$arr[0,9] = 2.375
$arr[0,10] = 2.8
$arr[0,11] = 3
$arr[0,12] = 2.375
$arr[0,13] = 2.8
$arr[0,14] = 3
$arr[0,15] = 0
$arr[0,16] = 1.75
$arr[0,17] = 2.05
$arr[0,18] = 1.75
$arr[0,19] = 1.8
$arr[0,20] = 2
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 2
$arr[0,24] = -1
$arr[0,25] = 1.05
$arr[0,26] = -1
$arr[0,27] = 1
$ws.Range($ws.Cells.Item(174,2), $ws.Cells.Item(174,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 7618724
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45283.35416666666
$arr[0,4] = 'Amanat Baghdad'
$arr[0,5] = 'Karbalaa FC'
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 'H'
$arr[0,9] = 2.875
$arr[0,10] = 2.75
$arr[0,11] = 2.5
$arr[0,12] = 2.875
$arr[0,13] = 2.75
$arr[0,14] = 2.5
$arr[0,15] = 0
$arr[0,16] = 2
$arr[0,17] = 1.8
$arr[0,18] = 2
$arr[0,19] = 2
$arr[0,20] = 1.8
$arr[0,21] = 1.875
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 1
$arr[0,25] = -1
$arr[0,26] = 1
$arr[0,27] = -1
$ws.Range($ws.Cells.Item(175,2), $ws.Cells.Item(175,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 7645264
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45291.58333333334
$arr[0,4] = 'Al Zawraa'
$arr[0,5] = 'Al Naft SC'
$arr[0,6] = 1
$arr[0,7] = 1
$arr[0,8] = 'D'
$arr[0,9] = 1.727
$arr[0,10] = 3.1
$arr[0,11] = 4.75
$arr[0,12] = 1.533
$arr[0,13] = 3.5
$arr[0,14] = 6.5
$arr[0,15] = -1
$arr[0,16] = 1.975
$arr[0,17] = 1.825
$arr[0,18] = 2
$arr[0,19] = 1.95
$arr[0,20] = 1.85
$arr[0,21] = -1
$arr[0,22] = 2.5
$arr[0,23] = -1
$arr[0,24] = -1
$arr[0,25] = 0.825
$arr[0,26] = 0
$arr[0,27] = 0
$ws.Range($ws.Cells.Item(186,2), $ws.Cells.Item(186,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 7648519
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45291.58333333334
$arr[0,4] = 'Al Najaf'
$arr[0,5] = 'Al Quwa Al Jawiya'
$arr[0,6] = 0
$arr[0,7] = 1
$arr[0,8] = 'A'
$arr[0,9] = 3.6
$arr[0,10] = 3.4
$arr[0,11] = 1.833
$arr[0,12] = 5.25
$arr[0,13] = 3.75
$arr[0,14] = 1.533
$arr[0,15] = 1
$arr[0,16] = 1.8
$arr[0,17] = 2
$arr[0,18] = 2
$arr[0,19] = 1.95
$arr[0,20] = 1.85
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.5329999999999999
$arr[0,24] = 0
$arr[0,25] = 0
$arr[0,26] = -1
$arr[0,27] = 0.8500000000000001
$ws.Range($ws.Cells.Item(187,2), $ws.Cells.Item(187,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 7811883
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45333.33333333334
$arr[0,4] = 'Al Karkh'
$arr[0,5] = 'Naft Maysan'
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 'D'
$arr[0,9] = 4.2
$arr[0,10] = 3.4
$arr[0,11] = 1.727
$arr[0,12] = 3.6
$arr[0,13] = 3.3
$arr[0,14] = 1.909
$arr[0,15] = 0.5
$arr[0,16] = 1.825
$arr[0,17] = 1.975
$arr[0,18] = 1.75
$arr[0,19] = 1.725
$arr[0,20] = 1.975
$arr[0,21] = -1
$arr[0,22] = 2.3
$arr[0,23] = -1
$arr[0,24] = 0.825
$arr[0,25] = -1
$arr[0,26] = -1
$arr[0,27] = 0.9750000000000001
$ws.Range($ws.Cells.Item(191,2), $ws.Cells.Item(191,29)).Value = $arr

$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 7811882
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45333.33333333334
$arr[0,4] = 'Amanat Baghdad'
$arr[0,5] = 'Naft AlWasat'
$arr[0,6] = 1
$arr[0,7] = 1
$arr[0,8] = 'D'
$arr[0,9] = 2.1
$arr[0,10] = 2.9
$arr[0,11] = 3.4
$arr[0,12] = 2.05
$arr[0,13] = 2.9
$arr[0,14] = 3.6
$arr[0,15] = -0.25
$arr[0,16] = 1.775
$arr[0,17] = 2.025
$arr[0,18] = 1.75
$arr[0,19] = 1.85
$arr[0,20] = 1.95
$arr[0,21] = -1
$arr[0,22] = 1.9
$arr[0,23] = -1
$arr[0,24] = -0.5
$arr[0,25] = 0.5125
$arr[0,26] = 0.425
$arr[0,27] = -0.5
$ws.Range($ws.Cells.Item(192,2), $ws.Cells.Item(192,29)).Value = $arr
# --- Row 211: update existing row in place (style already correct) ---
$arr = New-Object 'object[,]' 1,28
$arr[0,0] = 7870837
$arr[0,1] = 'Iraq League'
$arr[0,2] = 'Iraq League'
$arr[0,3] = 45345.54166666666
$arr[0,4] = 'Al Quwa Al Jawiya'
$arr[0,5] = 'Al Qasim SC'
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 'D'
$arr[0,9] = 1.4
$arr[0,10] = 3.75
$arr[0,11] = 7.5
$arr[0,12] = 1.333
$arr[0,13] = 3.75
$arr[0,14] = 10
$arr[0,15] = -1.25
$arr[0,16] = 1.8
$arr[0,17] = 2
$arr[0,18] = 2.5
$arr[0,19] = 1.975
$arr[0,20] = 1.825
$arr[0,21] = -1
$arr[0,22] = 2.75
$arr[0,23] = -1
$arr[0,24] = -1
$arr[0,25] = 1
$arr[0,26] = -1
$arr[0,27] = 0.825
$ws.Range($ws.Cells.Item(211,2), $ws.Cells.Item(211,29)).Value = $arr

# --- New rows 212-214: copy style from row 211 for column A (id) and E (date), then set values ---
$ws.Cells.Item(211, 1).Copy() | Out-Null
$ws.Cells.Item(212, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(211, 5).Copy() | Out-Null
$ws.Cells.Item(212, 5).PasteSpecial(-4122) | Out-Null
$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 210
$arr[0,1] = 7875109
$arr[0,2] = 'Iraq League'
$arr[0,3] = 'Iraq League'
$arr[0,4] = 45346.33333333334
$arr[0,5] = 'Al Naft SC'
$arr[0,6] = 'Al Minaa'
$arr[0,7] = 1
$arr[0,8] = 1
$arr[0,9] = 'D'
$arr[0,10] = 1.666
$arr[0,11] = 3.25
$arr[0,12] = 4.75
$arr[0,13] = 1.8
$arr[0,14] = 3
$arr[0,15] = 4.5
$arr[0,16] = -0.5
$arr[0,17] = 1.825
$arr[0,18] = 1.975
$arr[0,19] = 2
$arr[0,20] = 2.025
$arr[0,21] = 1.775
$arr[0,22] = -1
$arr[0,23] = 2
$arr[0,24] = -1
$arr[0,25] = -1
$arr[0,26] = 0.9750000000000001
$arr[0,27] = 0
$arr[0,28] = 0
$ws.Range($ws.Cells.Item(212,1), $ws.Cells.Item(212,29)).Value = $arr

$ws.Cells.Item(211, 1).Copy() | Out-Null
$ws.Cells.Item(213, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(211, 5).Copy() | Out-Null
$ws.Cells.Item(213, 5).PasteSpecial(-4122) | Out-Null
$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 211
$arr[0,1] = 7875108
$arr[0,2] = 'Iraq League'
$arr[0,3] = 'Iraq League'
$arr[0,4] = 45346.4375
$arr[0,5] = 'Al Zawraa'
$arr[0,6] = 'Al Karkh'
$arr[0,7] = 1
$arr[0,8] = 1
$arr[0,9] = 'D'
$arr[0,10] = 1.615
$arr[0,11] = 3.25
$arr[0,12] = 5.25
$arr[0,13] = 1.666
$arr[0,14] = 3.2
$arr[0,15] = 5
$arr[0,16] = -0.75
$arr[0,17] = 1.95
$arr[0,18] = 1.85
$arr[0,19] = 1.75
$arr[0,20] = 1.925
$arr[0,21] = 1.875
$arr[0,22] = -1
$arr[0,23] = 2.2
$arr[0,24] = -1
$arr[0,25] = -1
$arr[0,26] = 0.8500000000000001
$arr[0,27] = 0.4625
$arr[0,28] = -0.5
$ws.Range($ws.Cells.Item(213,1), $ws.Cells.Item(213,29)).Value = $arr

$ws.Cells.Item(211, 1).Copy() | Out-Null
$ws.Cells.Item(214, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(211, 5).Copy() | Out-Null
$ws.Cells.Item(214, 5).PasteSpecial(-4122) | Out-Null
$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 212
$arr[0,1] = 7864432
$arr[0,2] = 'Iraq League'
$arr[0,3] = 'Iraq League'
$arr[0,4] = 45346.54166666666
$arr[0,5] = 'Al Shorta SC'
$arr[0,6] = 'Al Kahrabaa'
$arr[0,7] = 2
$arr[0,8] = 1
$arr[0,9] = 'H'
$arr[0,10] = 1.666
$arr[0,11] = 3.25
$arr[0,12] = 4.75
$arr[0,13] = 1.6
$arr[0,14] = 3.3
$arr[0,15] = 5.75
$arr[0,16] = -0.75
$arr[0,17] = 1.775
$arr[0,18] = 2.025
$arr[0,19] = 2.25
$arr[0,20] = 1.875
$arr[0,21] = 1.925
$arr[0,22] = 0.6000000000000001
$arr[0,23] = -1
$arr[0,24] = -1
$arr[0,25] = 0.3875
$arr[0,26] = -0.5
$arr[0,27] = 0.875
$arr[0,28] = -1
$ws.Range($ws.Cells.Item(214,1), $ws.Cells.Item(214,29)).Value = $arr
